$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Output ("D38 NumberFormat: [" + $ws.Range("D38").NumberFormat + "]")
